$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string text: collapse the CRLF-wrapped labels onto a
#     single line using literal "<br>" tags (matches the new HTML-based
#     plot label rendering). Re-assign to every cell that shares the text
#     so no stale duplicate is left behind.
$global = "**Global**:<br>Implemented by<br>All other countries"
$highIncome = "**High-income**:<br>All other HICs and<br>not some MICs (such as China)"
$international = "**International**:<br>Some countries (e.g. EU, UK, Brazil)<br>and not others (e.g. U.S., China)"

for ($r = 2; $r -le 13; $r++) {
    $ws.Range("E$r").Value = $global
}
for ($r = 14; $r -le 25; $r++) {
    $ws.Range("E$r").Value = $highIncome
}
for ($r = 26; $r -le 37; $r++) {
    $ws.Range("E$r").Value = $international
}

# --- Refreshed mean / CI_low / CI_high estimates (re-run with the laptop's
#     ggplot2 3.5.1 after removing the fast-RU crop) for the four affected
#     along/y combinations: All x Global, Russia x Global, All x High-income,
#     Russia x High-income.
$ws.Range("B2").Value = 73.8161485502559
$ws.Range("C2").Value = 72.4542861110033
$ws.Range("D2").Value = 75.1780109895084

$ws.Range("B12").Value = 76.7816309654674
$ws.Range("C12").Value = 72.3067698657709
$ws.Range("D12").Value = 81.2564920651638

$ws.Range("B14").Value = 69.211130206042
$ws.Range("C14").Value = 67.7830711918507
$ws.Range("D14").Value = 70.6391892202333

$ws.Range("B24").Value = 69.4644695949361
$ws.Range("C24").Value = 64.5975635112594
$ws.Range("D24").Value = 74.3313756786127

"edits applied"
